$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date/number formatting from the previous data row (A2) onto the
# new row's date cell (A3) so it reuses the existing style index instead of
# Excel creating a brand new custom number format.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)  # xlPasteFormats

# Append the new data row (row 3) with the trading/sentiment sample values.
$ws.Range("A3").Value = 42605.648356481484
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 2133
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = "Named"
